$wb = $excel.ActiveWorkbook

# --- Overview sheet: Status text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" (shared by zh-cn / de-de rows) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B2").Value = "Handed back: in sync with en-US"
$zhcn.Range("B3").Value = "Handed back: in sync with en-US"

# New "Latest Target File" / "Latest Handback File" hyperlinked entries
$zhcn.Hyperlinks.Add($zhcn.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/0cba2f659c85a9055fdc72e5c4b28d1f84c83f48/e2e/9cc9d2d9-fc18-4f72-a0d6-d3d523bc7125.md", "", "", "9cc9d2d9-fc18-4f72-a0d6-d3d523bc7125.md")
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5f643c153db888935ca2415b7c6605297bfa132b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/9cc9d2d9-fc18-4f72-a0d6-d3d523bc7125.55a542bf3f62c3f5d38a4a45f875d0d00d75e8ca.zh-cn.xlf", "", "", "9cc9d2d9-fc18-4f72-a0d6-d3d523bc7125.55a542bf3f62c3f5d38a4a45f875d0d00d75e8ca.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/0cba2f659c85a9055fdc72e5c4b28d1f84c83f48/e2e/9cc9d2d9-fc18-4f72-a0d6-d3d523bc7125.md", "", "", "9cc9d2d9-fc18-4f72-a0d6-d3d523bc7125.md")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5f643c153db888935ca2415b7c6605297bfa132b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/9cc9d2d9-fc18-4f72-a0d6-d3d523bc7125.55a542bf3f62c3f5d38a4a45f875d0d00d75e8ca.zh-cn.xlf", "", "", "9cc9d2d9-fc18-4f72-a0d6-d3d523bc7125.55a542bf3f62c3f5d38a4a45f875d0d00d75e8ca.zh-cn.xlf")

$zhcn.Range("G2").Value = "2016-02-18 04:04:28"
$zhcn.Range("G3").Value = "2016-02-18 04:04:28"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B2").Value = "Handed back: in sync with en-US"
$dede.Range("B3").Value = "Handed back: in sync with en-US"

$dede.Hyperlinks.Add($dede.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/0cba2f659c85a9055fdc72e5c4b28d1f84c83f48/e2e/9cc9d2d9-fc18-4f72-a0d6-d3d523bc7125.md", "", "", "9cc9d2d9-fc18-4f72-a0d6-d3d523bc7125.md")
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/895eb9b29c83e3db4a02ee7ce76aee4e446931fd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/9cc9d2d9-fc18-4f72-a0d6-d3d523bc7125.55a542bf3f62c3f5d38a4a45f875d0d00d75e8ca.de-de.xlf", "", "", "9cc9d2d9-fc18-4f72-a0d6-d3d523bc7125.55a542bf3f62c3f5d38a4a45f875d0d00d75e8ca.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/0cba2f659c85a9055fdc72e5c4b28d1f84c83f48/e2e/9cc9d2d9-fc18-4f72-a0d6-d3d523bc7125.md", "", "", "9cc9d2d9-fc18-4f72-a0d6-d3d523bc7125.md")
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/895eb9b29c83e3db4a02ee7ce76aee4e446931fd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/9cc9d2d9-fc18-4f72-a0d6-d3d523bc7125.55a542bf3f62c3f5d38a4a45f875d0d00d75e8ca.de-de.xlf", "", "", "9cc9d2d9-fc18-4f72-a0d6-d3d523bc7125.55a542bf3f62c3f5d38a4a45f875d0d00d75e8ca.de-de.xlf")

$dede.Range("G2").Value = "2016-02-18 04:04:50"
$dede.Range("G3").Value = "2016-02-18 04:04:50"
